# Regenerate the "K" column (column G) of the save_data sheet.
# The previous values were computed from a "Strike#" based calculation;
# they are replaced here with the newly computed K values (s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new K value (column G), per the recalculated s_vals.
$kValues = @{
    2  = 5
    3  = 5
    4  = 1
    5  = 4
    6  = 4
    7  = 7
    8  = 1
    9  = 6
    10 = 8
    11 = 5
    12 = 5
    13 = 3
    14 = 7
    15 = 7
    16 = 7
    17 = 6
    18 = 9
    19 = 5
    20 = 8
    21 = 9
    22 = 5
    23 = 6
    24 = 10
    25 = 3
    26 = 2
    27 = 5
    28 = 6
    29 = 5
    30 = 6
    31 = 3
    32 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
